$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D updates are written as text (values contain dots used as
# thousands separators, e.g. "42.897.96", which Excel would otherwise try to
# parse/reformat as numbers). We temporarily force a text number format while
# assigning the values, then restore the original "Normal" style so no visible
# formatting change is left behind.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '42.897.96'
$ws.Range("E2").Value = '  +0.18%  '

$ws.Range("D3").Value = '2.304.70'

$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = '305.26'
$ws.Range("E5").Value = '  +2.19%  '

$ws.Range("D6").Value = '97.14'
$ws.Range("E6").Value = '  +0.43%  '

$ws.Range("E7").Value = '  -1.24%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").Value = '0.499'
$ws.Range("E9").Value = '  -0.81%  '

$ws.Range("D10").Value = '35.31'
$ws.Range("E10").Value = '  -0.80%  '

$ws.Range("D11").Value = '0.0786'
$ws.Range("E11").Value = '  -0.13%  '

$ws.Range("D12").Value = '18.71'
$ws.Range("E12").Value = '  +6.00%  '

$ws.Range("E13").Value = '  +1.98%  '

$ws.Range("E14").Value = '  +1.86%  '

$ws.Range("D15").Value = '2.663.22'

$ws.Range("D16").Value = '2.308.22'
$ws.Range("E16").Value = '  +0.57%  '

$ws.Range("D17").Value = '0.782'
$ws.Range("E17").Value = '  +1.04%  '

$ws.Range("D18").Value = '42.805.93'
$ws.Range("E18").Value = '  +0.14%  '

$ws.Range("D19").Value = '12.70'
$ws.Range("E19").Value = '  +1.36%  '

$ws.Range("E20").Value = '  -0.82%  '

$ws.Range("D21").Value = '6.03'
$ws.Range("E21").Value = '  -0.01%  '

$ws.Range("D22").Value = '67.25'
$ws.Range("E22").Value = '  -0.68%  '

$ws.Range("D23").Value = '236.34'
$ws.Range("E23").Value = '  -1.76%  '

$ws.Range("D24").Value = '2.15'
$ws.Range("E24").Value = '  +1.24%  '

$ws.Range("E25").Value = '  +0.12%  '

$ws.Range("E26").Value = '  +0.12%  '

$ws.Range("D27").Value = '24.79'
$ws.Range("E27").Value = '  -1.04%  '

$ws.Range("D28").Value = '166.67'
$ws.Range("E28").Value = '  +0.55%  '

$ws.Range("D29").Value = '2.06'
$ws.Range("E29").Value = '  +1.25%  '

$ws.Range("E30").Value = '  +0.24%  '

$ws.Range("D31").Value = '33.17'
$ws.Range("E31").Value = '  +1.35%  '

$ws.Range("E32").Value = '  +0.10%  '

$ws.Range("D33").Value = '18.19'
$ws.Range("E33").Value = '  +6.19%  '

$ws.Range("E34").Value = '  -0.17%  '

$ws.Range("D35").Value = '4.48'
$ws.Range("E35").Value = '  -5.68%  '

$ws.Range("D37").Value = '0.0688'
$ws.Range("E37").Value = '  +0.51%  '

$ws.Range("E38").Value = '  +0.71%  '

$ws.Range("E39").Value = '  -0.25%  '

$ws.Range("E40").Value = '  -0.52%  '

$ws.Range("E41").Value = '  -0.94%  '

$ws.Range("D42").Value = '1.999.55'
$ws.Range("E42").Value = '  -0.69%  '

$ws.Range("E43").Value = '  -0.47%  '

$ws.Range("D44").Value = '10.25'
$ws.Range("E44").Value = '  +1.89%  '

$ws.Range("D45").Value = '18.02'
$ws.Range("E45").Value = '  +5.20%  '

$ws.Range("E46").Value = '  +1.46%  '

$ws.Range("D47").Value = '2.78'
$ws.Range("E47").Value = '  +0.55%  '

$ws.Range("B48").Value = 'HuobiToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D48").Value = '2.90'
$ws.Range("E48").Value = '  +4.30%  '

$ws.Range("B49").Value = 'MultiversX'
$ws.Range("C49").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D49").Value = '53.66'
$ws.Range("E49").Value = '  +1.34%  '

$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.531.23'
$ws.Range("E50").Value = '  +0.70%  '

$ws.Range("D51").Value = '71.20'
$ws.Range("E51").Value = '  -0.76%  '

# Restore the default cell style for column D so no stray number formatting
# remains applied (matches original workbook formatting).
$dRange.Style = "Normal"
